$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("D3").Value = 221
$ws.Range("D42").Value = 1145

# Update view / selection state
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D2:D3").Select()
